# Weekly update: insert a new daily price record for "Espinaca" at
# Mercado Mayorista Lo Valledor de Santiago, pushing the existing
# records in rows 548:568 down by one row (to 549:569).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 548; Excel shifts rows 548:568 down to 549:569
# and carries the row formatting (incl. the date style in column D)
# down with them.
$ws.Rows.Item(548).Insert()

# Populate the newly inserted row with the latest weekly record.
$ws.Range("A548").Value = 6
$ws.Range("B548").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C548").Value = "Metropolitana"
$ws.Range("D548").Value = 44747
$ws.Range("E548").Value = 13
$ws.Range("F548").Value = 100112012
$ws.Range("G548").Value = "Espinaca"
$ws.Range("H548").Value = "Sin especificar"
$ws.Range("I548").Value = "Primera"
$ws.Range("J548").Value = 460
$ws.Range("K548").Value = 7500
$ws.Range("L548").Value = 8000
$ws.Range("M548").Value = 7707
$ws.Range("N548").Value = "`$/cuna 10 kilos"
$ws.Range("O548").Value = "Región Metropolitana"
$ws.Range("P548").Value = 771
$ws.Range("Q548").Value = 10
$ws.Range("R548").Value = "Hortaliza"
